$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 293, pushing existing row 293 (and all below it) down to 294.
$ws.Rows("293:293").Insert()

# Populate the newly inserted row 293 with the new weekly data point.
# Most fields mirror the row that used to be at 293 (now shifted to 294),
# except the date (D) and volume (M), which carry the new values.
$ws.Range("A293").Value = 4
$ws.Range("B293").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C293").Value = 'Los Lagos'
$ws.Range("D293").Value = 45015
$ws.Range("E293").Value = 10
$ws.Range("F293").Value = 'Fruta'
$ws.Range("G293").Value = 100108
$ws.Range("H293").Value = 'Tropicales y subtropicales'
$ws.Range("I293").Value = 100108002
$ws.Range("J293").Value = 'Mango'
$ws.Range("K293").Value = 'Sin especificar'
$ws.Range("L293").Value = 'Primera'
$ws.Range("M293").Value = 100
$ws.Range("N293").Value = 8000
$ws.Range("O293").Value = 8500
$ws.Range("P293").Value = 8250
$ws.Range("Q293").Value = '$/bandeja 4 kilos'
$ws.Range("R293").Value = 'Perú'
$ws.Range("S293").Value = 2062
$ws.Range("T293").Value = 4
